# The authored change (see commit diff) adds alt-text / description
# "QuizAnswer" to the "Smiley Face 3" shape (shape id 4) that lives on
# the slide with sldId="559" (cId="2705371955"), which is slide index 11
# in the deck's Slides collection (ppt/slides/slide11.xml).
#
# (All the other hunks in the source diff - the pc:chgInfo / collaboration
# revision log pruning in ppt/changesInfos/changesInfo1.xml and the cached
# "datetimeFigureOut" field text on the slide master / layouts / notes
# master - are PowerPoint's own automatic bookkeeping side effects of
# saving the deck on the day of the edit; they are not reachable/desirable
# to hand-author through the Shape/TextRange object model, since writing
# to a field-backed TextRange breaks the <a:fld> into a plain run instead
# of refreshing its cached text.)

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(11)
$shape = $slide.Shapes.Item(3)

Write-Host ("Target shape: " + $shape.Name + " (id=" + $shape.Id + ")")

$shape.AlternativeText = "QuizAnswer"

Write-Host ("New AlternativeText: " + $shape.AlternativeText)
